$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 270
$ws1.Range("F4").Value = 1906
$ws1.Range("F5").Value = 1576
$ws1.Range("F8").Value = 567
$ws1.Range("F9").Value = 132

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 270
$ws4.Range("F4").Value = 1906
$ws4.Range("F5").Value = 1576
$ws4.Range("F9").Value = 567
$ws4.Range("F10").Value = 132
